# Feature_Tracker.xlsx update
# Adds 8 new feature-request rows (30-37) to Sheet1:
#   - New icon / history related features requested by Andrew Knowles and Adim
#   - A handful of other feature requests from Weston Fiala
# None of these have a "Completed Version" (column C) yet, so they remain
# visible under the existing auto-filter (which hides rows where column C
# is not blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30 - Create Saved Roll from history
$ws.Range("A30").Value = "Create Saved Roll from history"
$ws.Range("B30").Value = "I have a roll in my history that I now want to save, add an option to save it."
$ws.Range("D30").Value = "Andrew Knowles"

# Row 31 - Re-Roll from history
$ws.Range("A31").Value = "Re-Roll from history"
$ws.Range("B31").Value = "I previously rolled some dice, and I want to be able to roll them again from the history tab"
$ws.Range("D31").Value = "Andrew Knowles"

# Row 32 - Separate roll sounds from crit sounds
$ws.Range("A32").Value = "Separate roll sounds from crit sounds"
$ws.Range("B32").Value = "I want the crit sounds to go off, but I do not want to have the roll sounds."
$ws.Range("D32").Value = "Adim - miguellicauco@gmail.com"

# Row 33 - Individual roll sounds
$ws.Range("A33").Value = "Individual roll sounds"
$ws.Range("B33").Value = "I want to be able to have a sound play when I roll an arbitrary number on an arbitrary dice"
$ws.Range("D33").Value = "Weston Fiala"

# Row 34 - Highlight min/max rolls
$ws.Range("A34").Value = "Highlight min/max rolls"
$ws.Range("B34").Value = "When I roll min/max I want it to be more visible in the roll display."
$ws.Range("D34").Value = "Weston Fiala"

# Row 35 - Better icons for saved rolls
$ws.Range("A35").Value = "Better icons for saved rolls"
$ws.Range("B35").Value = "The temp icons are not great. Find some better ones."
$ws.Range("D35").Value = "Weston Fiala"

# Row 36 - Dice with named faces
$ws.Range("A36").Value = "Dice with named faces"
$ws.Range("B36").Value = "I want to have dice with named faces instead of numbered faces"
$ws.Range("D36").Value = "Adim - miguellicauco@gmail.com"

# Row 37 - Don't lose all rolls on uninstall
$ws.Range("A37").Value = "Don't lose all rolls on uninstall"
$ws.Range("B37").Value = "I want to be able to uninstall and reinstall without losing my rolls. Screw you different development machines!"
$ws.Range("D37").Value = "Weston Fiala"

# Move the active selection past the newly-entered data, matching the
# author's cursor position after typing the last row.
$ws.Range("B38").Select()
